$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update vm_pu values for the 380 kV case (Case_5_68)
# Columns: B,C,D,E,F,I,J,K,L,M (and N for row 2 only)

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.030688913991656
$ws.Range("D2").Value = 1.033777895261851
$ws.Range("E2").Value = 1.039068203584429
$ws.Range("F2").Value = 1.046966858797791
$ws.Range("I2").Value = 1.031885535942883
$ws.Range("J2").Value = 1.035828773694291
$ws.Range("K2").Value = 1.036579266983556
$ws.Range("L2").Value = 1.041854443255225
$ws.Range("M2").Value = 1.049730814449256
$ws.Range("N2").Value = 1.005712725503983

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.031741231226775
$ws.Range("D3").Value = 1.034563260479991
$ws.Range("E3").Value = 1.040038460422519
$ws.Range("F3").Value = 1.048106192718482
$ws.Range("I3").Value = 1.032080558505501
$ws.Range("J3").Value = 1.036522126137798
$ws.Range("K3").Value = 1.037173746211644
$ws.Range("L3").Value = 1.042634420362071
$ws.Range("M3").Value = 1.050681047445336

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.032421978830613
$ws.Range("D4").Value = 1.035070716045394
$ws.Range("E4").Value = 1.040666572577727
$ws.Range("F4").Value = 1.04884400588789
$ws.Range("I4").Value = 1.032204654150423
$ws.Range("J4").Value = 1.036970043100419
$ws.Range("K4").Value = 1.03755703308098
$ws.Range("L4").Value = 1.043138792993386
$ws.Range("M4").Value = 1.051295917536339

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.032708124154603
$ws.Range("D5").Value = 1.035283875325801
$ws.Range("E5").Value = 1.040930700325731
$ws.Range("F5").Value = 1.049154323265423
$ws.Range("I5").Value = 1.032256321492679
$ws.Range("J5").Value = 1.037158172758817
$ws.Range("K5").Value = 1.037717835779568
$ws.Range("L5").Value = 1.043350753296417
$ws.Range("M5").Value = 1.051554410195557

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.032756166817261
$ws.Range("D6").Value = 1.035319655422098
$ws.Range("E6").Value = 1.040975052597772
$ws.Range("F6").Value = 1.049206435169524
$ws.Range("I6").Value = 1.032264967188873
$ws.Range("J6").Value = 1.037189750308486
$ws.Range("K6").Value = 1.037744815829531
$ws.Range("L6").Value = 1.04338633778654
$ws.Range("M6").Value = 1.05159781232165

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.032425802479163
$ws.Range("D7").Value = 1.03507356497882
$ws.Range("E7").Value = 1.040670101593319
$ws.Range("F7").Value = 1.048848151810992
$ws.Range("I7").Value = 1.032205346506906
$ws.Range("J7").Value = 1.036972557584637
$ws.Range("K7").Value = 1.037559183036168
$ws.Range("L7").Value = 1.04314162552391
$ws.Range("M7").Value = 1.051299371523063

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.031044585274102
$ws.Range("D8").Value = 1.03404346348822
$ws.Range("E8").Value = 1.039396045977315
$ws.Range("F8").Value = 1.047351780276541
$ws.Range("I8").Value = 1.031951878664176
$ws.Range("J8").Value = 1.036063246080657
$ws.Range("K8").Value = 1.036780459657748
$ws.Range("L8").Value = 1.042118107277318
$ws.Range("M8").Value = 1.050051948869959

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.028609376345163
$ws.Range("D9").Value = 1.032222745322431
$ws.Range("E9").Value = 1.037153238831544
$ws.Range("F9").Value = 1.044719485927107
$ws.Range("I9").Value = 1.031489194396967
$ws.Range("J9").Value = 1.03445535414898
$ws.Range("K9").Value = 1.035397694426084
$ws.Range("L9").Value = 1.040312059179624
$ws.Range("M9").Value = 1.047853879540023

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.026984992996845
$ws.Range("D10").Value = 1.031005240476398
$ws.Range("E10").Value = 1.035659554479716
$ws.Range("F10").Value = 1.042967647965572
$ws.Range("I10").Value = 1.031169976694748
$ws.Range("J10").Value = 1.033379692733724
$ws.Range("K10").Value = 1.034468780174821
$ws.Range("L10").Value = 1.039106369866965
$ws.Range("M10").Value = 1.046388535283522

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.026281395930686
$ws.Range("D11").Value = 1.030477178504828
$ws.Range("E11").Value = 1.03501313509336
$ws.Range("F11").Value = 1.04220979844481
$ws.Range("I11").Value = 1.031029203932895
$ws.Range("J11").Value = 1.032913035308694
$ws.Range("K11").Value = 1.034064877559208
$ws.Range("L11").Value = 1.038583900980822
$ws.Range("M11").Value = 1.045754032877813

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.026020013776208
$ws.Range("D12").Value = 1.030280901869584
$ws.Range("E12").Value = 1.034773079641204
$ws.Range("F12").Value = 1.041928405461056
$ws.Range("I12").Value = 1.030976531931428
$ws.Range("J12").Value = 1.032739564413446
$ws.Range("K12").Value = 1.033914598514627
$ws.Range("L12").Value = 1.038389773031814
$ws.Range("M12").Value = 1.045518350466792

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.026076082723042
$ws.Range("D13").Value = 1.030323009790796
$ws.Range("E13").Value = 1.034824569929549
$ws.Range("F13").Value = 1.041988760415856
$ws.Range("I13").Value = 1.030987847569339
$ws.Range("J13").Value = 1.032776780552771
$ws.Range("K13").Value = 1.033946845262465
$ws.Range("L13").Value = 1.038431416839815
$ws.Range("M13").Value = 1.045568905145192

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.026259790715033
$ws.Range("D14").Value = 1.030460956884819
$ws.Range("E14").Value = 1.03499329094344
$ws.Range("F14").Value = 1.042186536250796
$ws.Range("I14").Value = 1.031024857854761
$ws.Range("J14").Value = 1.032898698874761
$ws.Range("K14").Value = 1.034052460578182
$ws.Range("L14").Value = 1.038567855522897
$ws.Range("M14").Value = 1.045734551296487

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.026372974633472
$ws.Range("D15").Value = 1.030545933309158
$ws.Range("E15").Value = 1.03509725261747
$ws.Range("F15").Value = 1.042308406517207
$ws.Range("I15").Value = 1.031047610406151
$ws.Range("J15").Value = 1.032973799080085
$ws.Range("K15").Value = 1.0341175003243
$ws.Range("L15").Value = 1.038651911981281
$ws.Range("M15").Value = 1.045836611356729

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.02703168355727
$ws.Range("D16").Value = 1.031040267827895
$ws.Range("E16").Value = 1.035702462704827
$ws.Range("F16").Value = 1.043017958832799
$ws.Range("I16").Value = 1.031179265631301
$ws.Range("J16").Value = 1.033410644521688
$ws.Range("K16").Value = 1.034495550559537
$ws.Range("L16").Value = 1.039141036029976
$ws.Range("M16").Value = 1.046430645132446

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.027444812725435
$ws.Range("D17").Value = 1.03135011704933
$ws.Range("E17").Value = 1.036082190286905
$ws.Range("F17").Value = 1.04346323157793
$ws.Range("I17").Value = 1.031261167136932
$ws.Range("J17").Value = 1.033684428137307
$ws.Range("K17").Value = 1.034732242795875
$ws.Range("L17").Value = 1.039447744194719
$ws.Range("M17").Value = 1.046803267066119

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.02768576198494
$ws.Range("D18").Value = 1.031530762467251
$ws.Range("E18").Value = 1.036303713288071
$ws.Range("F18").Value = 1.043723019955627
$ws.Range("I18").Value = 1.031308692901478
$ws.Range("J18").Value = 1.033844035719236
$ws.Range("K18").Value = 1.034870139586724
$ws.Range("L18").Value = 1.039626603531484
$ws.Range("M18").Value = 1.047020611275255

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.027767915742266
$ws.Range("D19").Value = 1.031592343558395
$ws.Range("E19").Value = 1.036379252728076
$ws.Range("F19").Value = 1.043811612675361
$ws.Range("I19").Value = 1.031324856224625
$ws.Range("J19").Value = 1.033898443231947
$ws.Range("K19").Value = 1.034917131378947
$ws.Range("L19").Value = 1.039687583439695
$ws.Range("M19").Value = 1.047094720062507

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.027400490154436
$ws.Range("D20").Value = 1.031316881885772
$ws.Range("E20").Value = 1.036041445560049
$ws.Range("F20").Value = 1.043415450961699
$ws.Range("I20").Value = 1.031252405323061
$ws.Range("J20").Value = 1.033655062617753
$ws.Range("K20").Value = 1.03470686469121
$ws.Range("L20").Value = 1.039414841277742
$ws.Range("M20").Value = 1.046763288245709

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.026205694257544
$ws.Range("D21").Value = 1.030420338513968
$ws.Range("E21").Value = 1.034943605333275
$ws.Range("F21").Value = 1.042128293271012
$ws.Range("I21").Value = 1.031013969811108
$ws.Range("J21").Value = 1.032862800654307
$ws.Range("K21").Value = 1.034021366444235
$ws.Range("L21").Value = 1.038527679380796
$ws.Range("M21").Value = 1.045685772633981

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.025454276586432
$ws.Range("D22").Value = 1.029855888940358
$ws.Range("E22").Value = 1.034253659098088
$ws.Range("F22").Value = 1.041319619699191
$ws.Range("I22").Value = 1.030861841721372
$ws.Range("J22").Value = 1.032363901432299
$ws.Range("K22").Value = 1.033588910762012
$ws.Range("L22").Value = 1.037969539819271
$ws.Range("M22").Value = 1.045008295911444

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.025852636656291
$ws.Range("D23").Value = 1.030155185962223
$ws.Range("E23").Value = 1.034619383305418
$ws.Range("F23").Value = 1.041748254796789
$ws.Range("I23").Value = 1.030942697488676
$ws.Range("J23").Value = 1.032628450591577
$ws.Range("K23").Value = 1.033818301602469
$ws.Range("L23").Value = 1.038265452937314
$ws.Range("M23").Value = 1.045367439186826

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.027420517684114
$ws.Range("D24").Value = 1.031331899686264
$ws.Range("E24").Value = 1.03605985624207
$ws.Range("F24").Value = 1.043437040754171
$ws.Range("I24").Value = 1.031256365169647
$ws.Range("J24").Value = 1.033668331896888
$ws.Range("K24").Value = 1.034718332464109
$ws.Range("L24").Value = 1.039429708807975
$ws.Range("M24").Value = 1.046781352953167

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.029239095401842
$ws.Range("D25").Value = 1.032694098159833
$ws.Range("E25").Value = 1.037732791200228
$ws.Range("F25").Value = 1.045399463707819
$ws.Range("I25").Value = 1.031610707880881
$ws.Range("J25").Value = 1.034871691558839
$ws.Range("K25").Value = 1.035756420508819
$ws.Range("L25").Value = 1.04077925863368
$ws.Range("M25").Value = 1.048422127087495

